$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$types = @(
    @("1", "Desktop"),
    @("2", "Laptop"),
    @("3", "All in one")
)

$row = 2
for ($i = 0; $i -lt 9; $i++) {
    foreach ($pair in $types) {
        $idCell = $ws.Cells.Item($row, 1)
        $idCell.NumberFormat = "@"
        $idCell.Value = $pair[0]
        $ws.Cells.Item($row, 2).Value = $pair[1]
        $row = $row + 1
    }
}
